# Updated cryptos list on Sat Nov 23 04:53:38 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) / Volume(1h) (column E) columns for the
# crypto ranking table with newly scraped figures, and fixes the
# ARBITRUM / dogwifhat rows (44-45), which had swapped places along with
# their Link/Price/Volume data.
#
# NOTE: column D/E cells hold numeric-looking text (e.g. "1.52",
# "  +7.87%  ") that must stay plain text, matching the original
# workbook's inlineStr cells. A leading apostrophe forces Excel to treat
# the assigned value as text instead of auto-converting it to a number
# or percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''98.382.22'
$ws.Range("E2").Value = '''  -0.53%  '
$ws.Range("D3").Value = '''3.340.90'
$ws.Range("E3").Value = '''  -1.82%  '
$ws.Range("E4").Value = '''  +0.05%  '
$ws.Range("D5").Value = '''261.85'
$ws.Range("E5").Value = '''  +0.34%  '
$ws.Range("D6").Value = '''641.66'
$ws.Range("E6").Value = '''  +0.58%  '
$ws.Range("D7").Value = '''1.52'
$ws.Range("E7").Value = '''  +7.87%  '
$ws.Range("D8").Value = '''0.450'
$ws.Range("E8").Value = '''  +13.03%  '
$ws.Range("E9").Value = '''  +20.09%  '
$ws.Range("E10").Value = '''  +0.04%  '
$ws.Range("D11").Value = '''3.335.49'
$ws.Range("E11").Value = '''  -1.86%  '
$ws.Range("D12").Value = '''43.94'
$ws.Range("E12").Value = '''  +20.50%  '
$ws.Range("E13").Value = '''  +3.11%  '
$ws.Range("D14").Value = '''0.0000269'
$ws.Range("E14").Value = '''  +7.10%  '
$ws.Range("D15").Value = '''98.181.46'
$ws.Range("E15").Value = '''  -0.42%  '
$ws.Range("D16").Value = '''3.960.78'
$ws.Range("E16").Value = '''  -0.65%  '
$ws.Range("D17").Value = '''5.53'
$ws.Range("E17").Value = '''  -0.99%  '
$ws.Range("D18").Value = '''3.334.56'
$ws.Range("E18").Value = '''  -2.16%  '
$ws.Range("D19").Value = '''7.35'
$ws.Range("E19").Value = '''  +17.94%  '
$ws.Range("D20").Value = '''16.61'
$ws.Range("E20").Value = '''  +8.54%  '
$ws.Range("D21").Value = '''529.25'
$ws.Range("E21").Value = '''  +6.71%  '
$ws.Range("D22").Value = '''3.53'
$ws.Range("E22").Value = '''  -2.84%  '
$ws.Range("D23").Value = '''10.01'
$ws.Range("E23").Value = '''  +5.94%  '
$ws.Range("D24").Value = '''0.0000209'
$ws.Range("E24").Value = '''  -2.28%  '
$ws.Range("D25").Value = '''0.423'
$ws.Range("E25").Value = '''  +47.41%  '
$ws.Range("D26").Value = '''101.82'
$ws.Range("E26").Value = '''  +13.93%  '
$ws.Range("D27").Value = '''6.04'
$ws.Range("E27").Value = '''  +4.70%  '
$ws.Range("D28").Value = '''12.68'
$ws.Range("E28").Value = '''  +4.24%  '
$ws.Range("D29").Value = '''3.521.10'
$ws.Range("E29").Value = '''  +0.06%  '
$ws.Range("D30").Value = '''0.148'
$ws.Range("E30").Value = '''  +11.32%  '
$ws.Range("D31").Value = '''0.999'
$ws.Range("D32").Value = '''10.80'
$ws.Range("E32").Value = '''  +11.64%  '
$ws.Range("D33").Value = '''0.189'
$ws.Range("E33").Value = '''  -3.52%  '
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '''  +0.59%  '
$ws.Range("D35").Value = '''28.85'
$ws.Range("E35").Value = '''  +2.77%  '
$ws.Range("D36").Value = '''0.510'
$ws.Range("E36").Value = '''  +7.75%  '
$ws.Range("D37").Value = '''7.71'
$ws.Range("E37").Value = '''  +3.79%  '
$ws.Range("D38").Value = '''0.155'
$ws.Range("E38").Value = '''  +2.54%  '
$ws.Range("D39").Value = '''2.04'
$ws.Range("E39").Value = '''  +1.83%  '
$ws.Range("D40").Value = '''520.02'
$ws.Range("E40").Value = '''  +2.30%  '
$ws.Range("D41").Value = '''24.69'
$ws.Range("E41").Value = '''  -0.71%  '
$ws.Range("D42").Value = '''3.85'
$ws.Range("E42").Value = '''  +2.10%  '
$ws.Range("E43").Value = '''  +1.90%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '''0.802'
$ws.Range("E44").Value = '''  +1.56%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''3.29'
$ws.Range("E45").Value = '''  -3.67%  '
$ws.Range("E46").Value = '''  +0.00%  '
$ws.Range("D47").Value = '''0.0385'
$ws.Range("E47").Value = '''  +17.77%  '
$ws.Range("D48").Value = '''163.95'
$ws.Range("E48").Value = '''  +2.25%  '
$ws.Range("D49").Value = '''2.00'
$ws.Range("E49").Value = '''  +2.22%  '
$ws.Range("D50").Value = '''7.66'
$ws.Range("E50").Value = '''  +15.94%  '
$ws.Range("D51").Value = '''49.48'
$ws.Range("E51").Value = '''  +5.85%  '
